$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated symbol list values: Price (D), Volume 1h % (E), and Hora/Hour (G)
# for each coin row, per the upstream GitHub Actions data refresh.
$updates = @(
    @{Cell="D2"; Value='301.64'},
    @{Cell="E2"; Value='0.74%'},
    @{Cell="G2"; Value='23'},
    @{Cell="D3"; Value='32.13'},
    @{Cell="E3"; Value='1.35%'},
    @{Cell="G3"; Value='23'},
    @{Cell="D4"; Value='4.965'},
    @{Cell="E4"; Value='-3.58%'},
    @{Cell="G4"; Value='23'},
    @{Cell="D5"; Value='0.07885'},
    @{Cell="E5"; Value='-2.17%'},
    @{Cell="G5"; Value='23'},
    @{Cell="D6"; Value='2.119'},
    @{Cell="E6"; Value='-15.99%'},
    @{Cell="G6"; Value='23'},
    @{Cell="D7"; Value='7.806'},
    @{Cell="E7"; Value='0.28%'},
    @{Cell="G7"; Value='23'},
    @{Cell="D8"; Value='3.842'},
    @{Cell="E8"; Value='-1.90%'},
    @{Cell="G8"; Value='23'},
    @{Cell="D9"; Value='0.9272'},
    @{Cell="E9"; Value='-0.01%'},
    @{Cell="G9"; Value='23'},
    @{Cell="D10"; Value='0.1747'},
    @{Cell="E10"; Value='-0.82%'},
    @{Cell="G10"; Value='23'},
    @{Cell="D11"; Value='0.07908'},
    @{Cell="E11"; Value='7.59%'},
    @{Cell="G11"; Value='23'},
    @{Cell="D12"; Value='0.08626'},
    @{Cell="E12"; Value='-2.50%'},
    @{Cell="G12"; Value='23'},
    @{Cell="D13"; Value='0.03105'},
    @{Cell="E13"; Value='2.56%'},
    @{Cell="G13"; Value='23'},
    @{Cell="D14"; Value='0.1002'},
    @{Cell="E14"; Value='0.00%'},
    @{Cell="G14"; Value='23'},
    @{Cell="D15"; Value='0.001512'},
    @{Cell="E15"; Value='-1.30%'},
    @{Cell="G15"; Value='23'},
    @{Cell="D16"; Value='0.005914'},
    @{Cell="E16"; Value='3.22%'},
    @{Cell="G16"; Value='23'},
    @{Cell="E17"; Value='2,096.21%'},
    @{Cell="G17"; Value='23'},
    @{Cell="D18"; Value='3.461'},
    @{Cell="E18"; Value='-2.61%'},
    @{Cell="G18"; Value='23'},
    @{Cell="D19"; Value='2.256'},
    @{Cell="E19"; Value='-1.38%'},
    @{Cell="G19"; Value='23'},
    @{Cell="D20"; Value='0.3288'},
    @{Cell="E20"; Value='0.49%'},
    @{Cell="G20"; Value='23'},
    @{Cell="D21"; Value='0.1310'},
    @{Cell="E21"; Value='-2.31%'},
    @{Cell="G21"; Value='23'},
    @{Cell="D22"; Value='4.263'},
    @{Cell="E22"; Value='2.54%'},
    @{Cell="G22"; Value='23'},
    @{Cell="D23"; Value='0.1794'},
    @{Cell="E23"; Value='6.53%'},
    @{Cell="G23"; Value='23'},
    @{Cell="D24"; Value='0.04608'},
    @{Cell="E24"; Value='-0.45%'},
    @{Cell="G24"; Value='23'},
    @{Cell="D25"; Value='0.001237'},
    @{Cell="E25"; Value='-0.15%'},
    @{Cell="G25"; Value='23'},
    @{Cell="D26"; Value='0.004442'},
    @{Cell="E26"; Value='-1.82%'},
    @{Cell="G26"; Value='23'},
    @{Cell="D27"; Value='0.0001249'},
    @{Cell="E27"; Value='4.06%'},
    @{Cell="G27"; Value='23'},
    @{Cell="G28"; Value='23'},
    @{Cell="G29"; Value='23'},
    @{Cell="G30"; Value='23'},
    @{Cell="G31"; Value='23'},
    @{Cell="G32"; Value='23'},
    @{Cell="G33"; Value='23'},
    @{Cell="G34"; Value='23'},
    @{Cell="G35"; Value='23'},
    @{Cell="G36"; Value='23'},
    @{Cell="G37"; Value='23'},
    @{Cell="G38"; Value='23'},
    @{Cell="D39"; Value='0.01716'},
    @{Cell="E39"; Value='-2.16%'},
    @{Cell="G39"; Value='23'},
    @{Cell="D40"; Value='0.04778'},
    @{Cell="E40"; Value='3.79%'},
    @{Cell="G40"; Value='23'},
    @{Cell="D41"; Value='0.007431'},
    @{Cell="E41"; Value='7.62%'},
    @{Cell="G41"; Value='23'},
    @{Cell="D42"; Value='0.1358'},
    @{Cell="E42"; Value='-1.19%'},
    @{Cell="G42"; Value='23'},
    @{Cell="D43"; Value='0.002349'},
    @{Cell="E43"; Value='7.20%'},
    @{Cell="G43"; Value='23'},
    @{Cell="D44"; Value='0.01120'},
    @{Cell="E44"; Value='8.43%'},
    @{Cell="G44"; Value='23'},
    @{Cell="D45"; Value='0.00005981'},
    @{Cell="E45"; Value='-2.83%'},
    @{Cell="G45"; Value='23'},
    @{Cell="D46"; Value='0.00000000751'},
    @{Cell="E46"; Value='0.00%'},
    @{Cell="G46"; Value='23'},
    @{Cell="D47"; Value='0.003394'},
    @{Cell="E47"; Value='-59.66%'},
    @{Cell="G47"; Value='23'},
    @{Cell="D48"; Value='0.8205'},
    @{Cell="E48"; Value='9.62%'},
    @{Cell="G48"; Value='23'},
    @{Cell="D49"; Value='0.00002103'},
    @{Cell="E49"; Value='0.00%'},
    @{Cell="G49"; Value='23'},
    @{Cell="D50"; Value='0.0002003'},
    @{Cell="E50"; Value='0.00%'},
    @{Cell="G50"; Value='23'},
    @{Cell="G51"; Value='23'}
)

foreach ($item in $updates) {
    $cell = $ws.Range($item.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $item.Value
    $cell.Style = "Normal"
}
